$wb = $excel.ActiveWorkbook

# "임석렬" worksheet is the 4th (and active) sheet in this workbook.
$ws = $wb.Worksheets.Item(4)

# --- Row 10 used to be a blank template row; copy the formatting (styles +
# row height) from row 3, which already carries the exact style pattern
# (3,1,2,2,3,1) and the 75pt row height this row ends up with. ---
$ws.Range("A3:F3").Copy()
$ws.Range("A10:F10").PasteSpecial(-4122)
$ws.Rows.Item(10).RowHeight = $ws.Rows.Item(3).RowHeight

# NOTE: new shared strings are interned in first-use order, so the cells
# below are written in the same left-to-right-by-new-string order as the
# target sharedStrings.xml (A10, then F10, then B10) to land on shared
# string indices 98/99/100 respectively.

# --- A10: "OnlineController 작성" ---
$ws.Range("A10").Value = "OnlineController 작성"

# --- F10: "connection 객체 넘겨받는 부분 필요함" with the Korean words in 굴림 ---
$cellF = $ws.Range("F10")
$cellF.Value = "connection 객체 넘겨받는 부분 필요함"
$cellF.Characters(12, 2).Font.Name = "굴림"
$cellF.Characters(15, 4).Font.Name = "굴림"
$cellF.Characters(20, 2).Font.Name = "굴림"
$cellF.Characters(23, 3).Font.Name = "굴림"

# --- B10: "connection과 연동하여 작동하는 부분 작성 " ---
$cellB = $ws.Range("B10")
$cellB.Value = "connection과 연동하여 작동하는 부분 작성 "
$cellB.Characters(11, 18).Font.Name = "굴림"

# --- C10 / D10: assigned/completed dates ---
$ws.Range("C10").Value = 42153
$ws.Range("D10").Value = 42153

# --- E10: "작성함" (reuses the existing shared string) ---
$ws.Range("E10").Value = "작성함"

# --- Move the active selection from F10 to E11 ---
[void]$ws.Range("E11").Select()
